$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgabenliste Projekt 1")

# --- Row 13: Risikenliste anfertigen -> person change Jana -> Viktoria
$ws.Range("D13").Value = "Viktoria"

# --- Row 20: Textur - Boden (externe Quelle) -> person change Viktoria -> Tobias, actual completion date added
$ws.Range("D20").Value = "Tobias"
$ws.Range("H19").Copy($ws.Range("H20"))
$ws.Range("H20").Value = 42342

# --- Row 22: Holztexturen -> actual completion date added
$ws.Range("H21").Copy($ws.Range("H22"))
$ws.Range("H22").Value = 42341

# --- Row 29: (M) - Schrank + Inhalt(...) -> % erledigt 60% -> 100%
$ws.Range("E29").Value = 1

# --- Row 30: (M) - Bett + schlafenden Jungen -> % erledigt 0% -> 100%, actual completion date added
$ws.Range("E30").Value = 1
$ws.Range("H29").Copy($ws.Range("H30"))
$ws.Range("H30").Value = 42342

# --- Row 36: (S) - Schrank, Mobile, Fenster-Szene fertig stellen -> % erledigt 0% -> 100%, completion date added
$ws.Range("E36").Value = 1
$ws.Range("H35").Copy($ws.Range("H36"))
$ws.Range("H36").Value = 42343

# --- Row 37: (S) - Bett + Kommode + Junge Szene fertig stellen -> % erledigt 0% -> 100%, completion date added
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 1
$ws.Range("H35").Copy($ws.Range("H37"))
$ws.Range("H37").Value = 42343

# --- Row 39: Beleuchtung - Nacht -> completion date added
$ws.Range("H38").Copy($ws.Range("H39"))
$ws.Range("H39").Value = 42343

# --- Row 40: Beleuchtung - Morgensonne -> person Viktoria, % erledigt 0% -> 100%, completion date added
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = ""
$ws.Range("H38").Copy($ws.Range("H40"))
$ws.Range("H40").Value = 42343

# --- Row 41: (R) - Rendern der Beispielbilder -> completion date added
$ws.Range("H38").Copy($ws.Range("H41"))
$ws.Range("H41").Value = 42343

# --- Row 42: (A) - Animation der Shots 1-4 -> renamed
$ws.Range("B42").Value = "(A) - Animation der Shots 1-5 + 13"

# --- Row 43: (A) - Animation der Shots 5 - 9 -> renamed
$ws.Range("B43").Value = "(A) - Animation der Shots 6 - 10"

# --- Row 44: (A) - Animation der Shots 10 - 14 -> renamed
$ws.Range("B44").Value = "(A) - Animation der Shots 11 - 15, außer 13"

# --- Row 46: (A) - Erstellung von Anfang und Abspann -> due date changed
$ws.Range("C46").Value = 42350

# --- Row 52: Abgabe finales Video -> renamed
$ws.Range("B52").Value = "Abgabe finales Video + Präsentation"

# --- Move active selection to reflect latest edit location
$ws.Range("E41").Select()
